# Clear the "No URL" placeholder text from the four department rows that
# don't actually have a catalog URL (Industrial Arts and Vocational
# Education, Basic Pharmaceutical Sciences, Pharmacy Practice, Economic,
# Social and Administrative Pharmacy). These cells keep their existing
# style but lose their shared-string value, becoming blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()
